# Reposition/resize the two pictures on slide 2 (indices 2 and 3 in the
# shape collection: 1=Title, 2=Picture 2, 3=Picture 3, 4=TextBox 4).
#
# EMU -> points conversion: 1 pt = 12700 EMU

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Picture 2 (rId2): off 274320,914400 ext 5486400,4114800 -> off 1371600,2286000 ext 3657600,2743200
$pic1 = $s.Shapes.Item(2)
$pic1.Left   = 1371600 / 12700
$pic1.Top    = 2286000 / 12700
$pic1.Width  = 3657600 / 12700
$pic1.Height = 2743200 / 12700

# Picture 3 (rId3): off 6370323,914400 ext 5486400,4114800 -> off 6400800,2286000 ext 3657600,2743200
$pic2 = $s.Shapes.Item(3)
$pic2.Left   = 6400800 / 12700
$pic2.Top    = 2286000 / 12700
$pic2.Width  = 3657600 / 12700
$pic2.Height = 2743200 / 12700
